# Add a new "Directly send Outlook mail" check row to the KO Checklist
# workbook (sheet 1 = "워크플로우"), right after the existing last row (33).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clone formatting (column styles + row height) from the last existing
# data row (33) onto the new row (34) before filling in values, so the
# new row matches the look of the rest of the table.
$ws.Range("A33:G33").Copy()
$ws.Range("A34:G34").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(34).RowHeight = 68

# Fill in the new check's data.
$ws.Range("A34").Value2 = "No"
$ws.Range("B34").Value2 = "Directly send Outlook mail"
$ws.Range("C34").Value2 = "Checks\Custom\DirectlySendOutlookMail.xaml"
$ws.Range("D34").Value2 = ""
$ws.Range("E34").Value2 = "Fix"
$ws.Range("F34").Value2 = "According to the CoE (Centre of Excellence)'s security policies, robots should not be allowed to directly send emails. Instead, emails created by robots should be saved as drafts and then reviewed by humans before sending."
$ws.Range("G34").Value2 = "Check the IsDraft property of Send Outlook Mail Message activities, as specified by the CoE's security policies."
